$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: write all literal text values first (as apostrophe-prefixed strings
# so look-like-numbers stay text), while the new rows still carry the sheet
# default (General) formatting.
$ws.Range("D232").Value = "'16264"
$ws.Range("E232").Value = "'Bacteria:Binomial (genus species)"
$ws.Range("F232").Value = "'2: 851"
$ws.Range("G232").Value = "'2: 873"
$ws.Range("I232").Value = "'Acinetobacter baumannii"
$ws.Range("L232").Value = "'Sonia"
$ws.Range("M232").Value = "'11/8/18 14:38:00"
$ws.Range("H232").Value = 0
$ws.Range("J232").Value = 23
$ws.Range("K232").Value = 0.133877

$ws.Range("D233").Value = "'16715"
$ws.Range("E233").Value = "'Bacteria:Binomial (genus species)"
$ws.Range("F233").Value = "'1: 29"
$ws.Range("G233").Value = "'1: 49"
$ws.Range("I233").Value = "'Klebsiella pneumoniae"
$ws.Range("L233").Value = "'Sonia"
$ws.Range("M233").Value = "'11/8/18 14:38:00"
$ws.Range("H233").Value = 0
$ws.Range("J233").Value = 21
$ws.Range("K233").Value = 0.152683

$ws.Range("D234").Value = "'16264"
$ws.Range("E234").Value = "'Event month"
$ws.Range("F234").Value = "'2: 2613"
$ws.Range("G234").Value = "'2: 2617"
$ws.Range("I234").Value = "'April"
$ws.Range("L234").Value = "'Sonia"
$ws.Range("M234").Value = "'11/13/18 08:33:00"
$ws.Range("H234").Value = 0
$ws.Range("J234").Value = 5
$ws.Range("K234").Value = 0.029104

$ws.Range("D235").Value = "'16264"
$ws.Range("E235").Value = "'Event month"
$ws.Range("F235").Value = "'2: 2590"
$ws.Range("G235").Value = "'2: 2598"
$ws.Range("I235").Value = "'September"
$ws.Range("L235").Value = "'Sonia"
$ws.Range("M235").Value = "'11/13/18 08:33:00"
$ws.Range("H235").Value = 0
$ws.Range("J235").Value = 9
$ws.Range("K235").Value = 0.052386

$ws.Range("D236").Value = "'16264"
$ws.Range("E236").Value = "'Event year"
$ws.Range("F236").Value = "'2: 2600"
$ws.Range("G236").Value = "'2: 2603"
$ws.Range("I236").Value = "'2002"
$ws.Range("L236").Value = "'Sonia"
$ws.Range("M236").Value = "'11/13/18 08:33:00"
$ws.Range("H236").Value = 0
$ws.Range("J236").Value = 4
$ws.Range("K236").Value = 0.023283

$ws.Range("D237").Value = "'16264"
$ws.Range("E237").Value = "'Event year"
$ws.Range("F237").Value = "'2: 2620"
$ws.Range("G237").Value = "'2: 2623"
$ws.Range("I237").Value = "'2005"
$ws.Range("L237").Value = "'Sonia"
$ws.Range("M237").Value = "'11/13/18 08:33:00"
$ws.Range("H237").Value = 0
$ws.Range("J237").Value = 4
$ws.Range("K237").Value = 0.023283

$ws.Range("D238").Value = "'16264"
$ws.Range("E238").Value = "'B"
$ws.Range("F238").Value = "'2: 2620"
$ws.Range("G238").Value = "'2: 2623"
$ws.Range("I238").Value = "'2005"
$ws.Range("L238").Value = "'Sonia"
$ws.Range("M238").Value = "'11/13/18 08:33:00"
$ws.Range("H238").Value = 0
$ws.Range("J238").Value = 4
$ws.Range("K238").Value = 0.023283

$ws.Range("D239").Value = "'16264"
$ws.Range("E239").Value = "'B"
$ws.Range("F239").Value = "'2: 2613"
$ws.Range("G239").Value = "'2: 2617"
$ws.Range("I239").Value = "'April"
$ws.Range("L239").Value = "'Sonia"
$ws.Range("M239").Value = "'11/13/18 08:33:00"
$ws.Range("H239").Value = 0
$ws.Range("J239").Value = 5
$ws.Range("K239").Value = 0.029104

$ws.Range("D240").Value = "'16264"
$ws.Range("E240").Value = "'A"
$ws.Range("F240").Value = "'2: 2600"
$ws.Range("G240").Value = "'2: 2603"
$ws.Range("I240").Value = "'2002"
$ws.Range("L240").Value = "'Sonia"
$ws.Range("M240").Value = "'11/13/18 08:33:00"
$ws.Range("H240").Value = 0
$ws.Range("J240").Value = 4
$ws.Range("K240").Value = 0.023283

$ws.Range("D241").Value = "'16264"
$ws.Range("E241").Value = "'A"
$ws.Range("F241").Value = "'2: 2590"
$ws.Range("G241").Value = "'2: 2598"
$ws.Range("I241").Value = "'September"
$ws.Range("L241").Value = "'Sonia"
$ws.Range("M241").Value = "'11/13/18 08:33:00"
$ws.Range("H241").Value = 0
$ws.Range("J241").Value = 9
$ws.Range("K241").Value = 0.052386

$ws.Range("D242").Value = "'16726"
$ws.Range("E242").Value = "'Event year"
$ws.Range("F242").Value = "'1: 3157"
$ws.Range("G242").Value = "'1: 3160"
$ws.Range("I242").Value = "'2007"
$ws.Range("L242").Value = "'Sonia"
$ws.Range("M242").Value = "'11/13/18 08:35:00"
$ws.Range("H242").Value = 0
$ws.Range("J242").Value = 4
$ws.Range("K242").Value = 0.022967

$ws.Range("D243").Value = "'16726"
$ws.Range("E243").Value = "'Event year"
$ws.Range("F243").Value = "'1: 3165"
$ws.Range("G243").Value = "'1: 3168"
$ws.Range("I243").Value = "'2009"
$ws.Range("L243").Value = "'Sonia"
$ws.Range("M243").Value = "'11/13/18 08:35:00"
$ws.Range("H243").Value = 0
$ws.Range("J243").Value = 4
$ws.Range("K243").Value = 0.022967

# Step 2: stamp column A/B/C (constant across all rows) to match the rest of the table.
$ws.Range("A232:A243").Value = "●"

# Step 3: copy the formatting (styles only) from the last existing data row (231)
# down across the newly added rows, without touching the values just written.
$ws.Range("A231:M231").Copy()
$ws.Range("A232:M243").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Step 4: match row height (16 points, same as the rest of the data rows).
$ws.Range("A232:M243").RowHeight = 16

Write-Output "done"